$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Estimates")

# Mark stories as Completed (checks the "Completed" checkbox column E)
# Completed Points (F) / Completed Hours (G) recalc automatically from Table1 formulas.
$ws.Range("E30").Value = $true
$ws.Range("E31").Value = $true
$ws.Range("E34").Value = $true
$ws.Range("E48").Value = $true

# Un-mark story 76 as Completed
$ws.Range("E76").ClearContents()

# Rows 79 & 87 lost their one-off "applyNumberFormat" styling (D back to the
# standard centered style used elsewhere in the table; F/G back to default).
$ws.Range("D79").HorizontalAlignment = -4108
$ws.Range("F79").ClearFormats()
$ws.Range("G79").ClearFormats()

$ws.Range("D87").HorizontalAlignment = -4108
$ws.Range("F87").ClearFormats()
$ws.Range("G87").ClearFormats()

# Restore the view: scrolled back near the top, with E80 selected.
$null = $ws.Range("E80").Select()
